# Auto-generated edit script for 广州-漫展信息.xlsx
# Applies the bilibili-event-data refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

function Set-NumCell($ws, $cellRef, $num) {
    $ws.Range($cellRef).Value = $num
}

function Set-TextCell($ws, $cellRef, $text) {
    # Force text storage so date-like / numeric-like strings are not
    # auto-converted by Excel, then reset the style so no extra
    # formatting (e.g. quote-prefix) lingers on the cell.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# --- Sheet 1 ("展览") : refresh "想去人数" (want-to-go) counters ---
$ws1 = $wb.Worksheets.Item(1)
Set-NumCell $ws1 "F4" 799
Set-NumCell $ws1 "F6" 1013
Set-NumCell $ws1 "F7" 1078
Set-NumCell $ws1 "F9" 148
Set-NumCell $ws1 "F10" 491
Set-NumCell $ws1 "F11" 209
Set-NumCell $ws1 "F12" 42
Set-NumCell $ws1 "F13" 1185
Set-NumCell $ws1 "F14" 26315
Set-NumCell $ws1 "F15" 2868
Set-NumCell $ws1 "F16" 17
Set-NumCell $ws1 "F17" 201
Set-NumCell $ws1 "F18" 403
Set-NumCell $ws1 "F21" 486
Set-NumCell $ws1 "F23" 192
Set-NumCell $ws1 "F24" 305
Set-NumCell $ws1 "F27" 146
Set-NumCell $ws1 "F28" 61
Set-NumCell $ws1 "F29" 441
Set-NumCell $ws1 "F32" 532
Set-NumCell $ws1 "F33" 215

# --- Sheet 2 ("演出") : refresh "想去人数" (want-to-go) counters ---
$ws2 = $wb.Worksheets.Item(2)
Set-NumCell $ws2 "F6" 329
Set-NumCell $ws2 "F7" 547
Set-NumCell $ws2 "F10" 4177
Set-NumCell $ws2 "F21" 4188

# --- Sheet 3 ("本地生活") : refresh "想去人数" (want-to-go) counters ---
$ws3 = $wb.Worksheets.Item(3)
Set-NumCell $ws3 "F4" 1021

# --- Sheet 4 ("全部类型") : refresh counters + insert new event row ---
$ws4 = $wb.Worksheets.Item(4)
Set-NumCell $ws4 "F4" 1021
Set-NumCell $ws4 "F6" 799
Set-NumCell $ws4 "F11" 329
Set-NumCell $ws4 "F13" 548
Set-NumCell $ws4 "F14" 1013
Set-NumCell $ws4 "F15" 1078
Set-NumCell $ws4 "F16" 148
Set-NumCell $ws4 "F17" 491
Set-NumCell $ws4 "F18" 209
Set-NumCell $ws4 "F19" 42
Set-NumCell $ws4 "F20" 1185
Set-NumCell $ws4 "F28" 201
Set-NumCell $ws4 "F31" 403
Set-NumCell $ws4 "F35" 486
Set-NumCell $ws4 "F37" 305
Set-NumCell $ws4 "F41" 146
Set-NumCell $ws4 "F42" 61
Set-NumCell $ws4 "F47" 532
Set-NumCell $ws4 "F48" 215

# A new performance event ("广州·萃火虫动漫游戏嘉年华...")
# now sorts into row 21 of the combined sheet; rows 21-26 shift down one
# row each, and the former row 27 event drops off this combined view.
# row 21
Set-TextCell $ws4 "B21" "2024-07-19"
Set-TextCell $ws4 "C21" "广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园"
Set-TextCell $ws4 "D21" "新港东路1000号 保利世贸博览馆"
Set-TextCell $ws4 "E21" "2024.07.19 09:00-07.22 17:00"
Set-NumCell $ws4 "F21" 26315
Set-TextCell $ws4 "G21" "已售罄"
Set-TextCell $ws4 "H21" "https://show.bilibili.com/platform/detail.html?id=87210"
Set-TextCell $ws4 "I21" "//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg"

# row 22
Set-TextCell $ws4 "B22" "2024-07-20"
Set-TextCell $ws4 "C22" "广州·冰兔2024线下live「过去和未来」"
Set-TextCell $ws4 "D22" "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）"
Set-TextCell $ws4 "E22" "2024.07.20 20:00-07.20 22:00"
Set-NumCell $ws4 "F22" 64
Set-NumCell $ws4 "G22" 198
Set-TextCell $ws4 "H22" "https://show.bilibili.com/platform/detail.html?id=87546"
Set-TextCell $ws4 "I22" "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg"

# row 23
Set-TextCell $ws4 "B23" "2024-07-20"
Set-TextCell $ws4 "C23" "广州·跨越二次元ACG神级动漫世界巡回演唱会"
Set-TextCell $ws4 "D23" "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
Set-TextCell $ws4 "E23" "2024.07.20 19:30-07.20 21:10"
Set-NumCell $ws4 "F23" 254
Set-NumCell $ws4 "G23" 280
Set-TextCell $ws4 "H23" "https://show.bilibili.com/platform/detail.html?id=85353"
Set-TextCell $ws4 "I23" "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"

# row 24
Set-TextCell $ws4 "B24" "2024-07-21"
Set-TextCell $ws4 "C24" "广州·昨日重现——唯美英文经典歌曲演唱会"
Set-TextCell $ws4 "D24" "东风中路299号 广州中山纪念堂"
Set-TextCell $ws4 "E24" "2024.07.21 19:30-07.21 21:30"
Set-NumCell $ws4 "F24" 2
Set-NumCell $ws4 "G24" 100
Set-TextCell $ws4 "H24" "https://show.bilibili.com/platform/detail.html?id=86802"
Set-TextCell $ws4 "I24" "//i1.hdslb.com/bfs/openplatform/202405/DR8AvmXe1716802703006.jpeg"

# row 25
Set-TextCell $ws4 "B25" "2024-07-21"
Set-TextCell $ws4 "C25" "广州·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024"
Set-TextCell $ws4 "D25" "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
Set-TextCell $ws4 "E25" "2024.07.21 14:30-07.21 16:00"
Set-NumCell $ws4 "F25" 174
Set-NumCell $ws4 "G25" 280
Set-TextCell $ws4 "H25" "https://show.bilibili.com/platform/detail.html?id=87034"
Set-TextCell $ws4 "I25" "//i1.hdslb.com/bfs/openplatform/202406/LINsP2ui1717741701901.png"

# row 26
Set-TextCell $ws4 "B26" "2024-07-26"
Set-TextCell $ws4 "C26" "广州·【早鸟8折】“浪漫古典Ⅱ”百年经典传世名曲烛光音乐会 "
Set-TextCell $ws4 "D26" "广州市二沙岛晴波路33号  星海音乐厅（交响乐演奏厅）"
Set-TextCell $ws4 "E26" "2024.07.26 20:00-07.26 21:30"
Set-NumCell $ws4 "F26" 1
Set-NumCell $ws4 "G26" 144
Set-TextCell $ws4 "H26" "https://show.bilibili.com/platform/detail.html?id=87726"
Set-TextCell $ws4 "I26" "//i0.hdslb.com/bfs/openplatform/202406/A8vhVlhn1717575084179.png"

# row 27
Set-TextCell $ws4 "B27" "2024-07-26"
Set-TextCell $ws4 "C27" "广州·萨克斯王子安德鲁·杨——2024经典&流行音乐巡回演出"
Set-TextCell $ws4 "D27" "龙凤街道革新路124号太古仓码头5号仓 广州太空间Live House"
Set-TextCell $ws4 "E27" "2024.07.26 20:00-07.26 21:30"
Set-NumCell $ws4 "F27" 3
Set-NumCell $ws4 "G27" 280
Set-TextCell $ws4 "H27" "https://show.bilibili.com/platform/detail.html?id=86635"
Set-TextCell $ws4 "I27" "//i1.hdslb.com/bfs/openplatform/202405/rciNih361716802006584.jpeg"

Write-Output "edit complete"
